# Regenerate merged AHB files
#
# - Rename header columns from the "*_old" / "*_new" naming scheme to the
#   "*_FV2404" / "*_FV2410" naming scheme (columns A:J and L:U of row 1;
#   column K, "diff", is left untouched).
# - Turn the data range A1:U92 into a proper Excel Table ("Table1") with an
#   AutoFilter.
# - Freeze the header row so it stays visible while scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename header row cells (row 1).
$fv2404Headers = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

$fv2410Headers = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)

# Columns A (1) .. J (10)
for ($i = 0; $i -lt $fv2404Headers.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $fv2404Headers[$i]
}

# Column K (11) = "diff" stays as-is.

# Columns L (12) .. U (21)
for ($i = 0; $i -lt $fv2410Headers.Count; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $fv2410Headers[$i]
}

# 2) Convert the used range A1:U92 into an Excel Table ("Table1") with an
#    AutoFilter on the header row.
$tableRange = $ws.Range("A1:U92")
$table = $ws.ListObjects.Add(1, $tableRange, 0, 1)
$table.Name = "Table1"

# 3) Freeze panes below row 1 (header row stays visible while scrolling).
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
